$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.980.91'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.242.13'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.82%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.634'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.28%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '69.81'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.84%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.557'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0987'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '58.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '35.45'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.87%  '
$ws.Range('E13').Value = '  -2.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.80'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.575.23'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.862'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -5.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.245.95'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.953.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0978'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.36'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.47'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.35%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.64'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.87%  '
$ws.Range('E29').Value = '  -8.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '168.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.62'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.32%  '
$ws.Range('E32').Value = '  -7.44%  '
$ws.Range('E33').Value = '  -7.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.45'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0720'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '21.98'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.19%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.25'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.75%  '
$ws.Range('B40').Value = 'THORChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.57%  '
$ws.Range('E41').Value = '  -4.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '66.93'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.07'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.34%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.93'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.100'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -10.23%  '
$ws.Range('E46').Value = '  -7.02%  '
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.36'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.43%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.94%  '
$ws.Range('B51').Value = 'Celestia'
$ws.Range('C51').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.97'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.27%  '
